# Saldo.xlsx update
# ------------------
# The workbook lists one account balance per row (Conta / Nome / Saldo).
# This script removes a batch of rows that dropped out of the latest export,
# updates one balance, and appends the accounts that are new in this export.
#
# Rows are located by searching column A (account number) rather than by a
# fixed row index, since every delete/insert shifts everything below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Find-AccountRow($sheet, $account) {
    $last = $sheet.UsedRange.Rows.Count
    for ($i = 1; $i -le $last; $i++) {
        $v = $sheet.Cells.Item($i, 1).Value()
        if ("$v" -eq $account) {
            return $i
        }
    }
    return -1
}

function Remove-AccountRow($sheet, $account) {
    $row = Find-AccountRow $sheet $account
    if ($row -gt 0) {
        $sheet.Rows.Item($row).Delete()
    }
}

function Set-AccountCell($sheet, $row, $account, $name, $saldo) {
    # Prefix with an apostrophe so Excel keeps the zero-padded account
    # number as text instead of inferring a number and dropping the
    # leading zeros.
    $sheet.Cells.Item($row, 1).Value = "'" + $account
    $sheet.Cells.Item($row, 2).Value = $name
    $sheet.Cells.Item($row, 3).Value = $saldo
}

# 1) Accounts removed entirely from this export.
Remove-AccountRow $ws "004361159"   # HFR            117415.34
Remove-AccountRow $ws "004497875"   # HENRIQUE        85242.82
Remove-AccountRow $ws "001882235"   # LAGO            75338.4 (old balance; re-added below with a new balance)
Remove-AccountRow $ws "004334062"   # MERG            49818.04
Remove-AccountRow $ws "004911541"   # TIAGO           47684.78
Remove-AccountRow $ws "004644524"   # PAULO           45425.49
Remove-AccountRow $ws "004328934"   # VALERIA         41880.4
Remove-AccountRow $ws "004363250"   # HELIO           39930.27
Remove-AccountRow $ws "004586209"   # ROBERIO         24806.36
Remove-AccountRow $ws "004453132"   # BRUNO           23888.05
Remove-AccountRow $ws "004385806"   # ANILSON         20674.33
Remove-AccountRow $ws "004479463"   # HENRIQUE        15939.47
Remove-AccountRow $ws "005121919"   # JORGE           15716.53
Remove-AccountRow $ws "004497825"   # PRISCILLA        2183.82

# 2) New account inserted right after SANDRA (004879567), where the two
#    removed HENRIQUE/LAGO rows used to be.
$afterRow = Find-AccountRow $ws "004879567"
$newRow = $afterRow + 1
$ws.Rows.Item($newRow).Insert()
Set-AccountCell $ws $newRow "005203796" "MARCIA" 74040.3

# 3) DANILO's balance changed from 9423.09 to 4000.
$danilo = Find-AccountRow $ws "005055226"
$ws.Cells.Item($danilo, 3).Value = 4000

# 4) LAGO re-appears further down the list with a new (much smaller) balance.
$afterLeone = Find-AccountRow $ws "004334158"   # LEONE, 994.66
$lagoRow = $afterLeone + 1
$ws.Rows.Item($lagoRow).Insert()
Set-AccountCell $ws $lagoRow "001882235" "LAGO" 966.83
